$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.286.48"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "2.619.27"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "522.65"
$ws.Range("E5").Value = "  +1.05%  "

$ws.Range("D6").Value = "144.24"
$ws.Range("E6").Value = "  +1.30%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("D9").Value = "2.618.19"
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").Value = "6.64"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").Value = "3.056.71"
$ws.Range("E14").Value = "  -0.18%  "

$ws.Range("D15").Value = "58.389.62"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "20.53"
$ws.Range("E16").Value = "  -2.04%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.641.61"
$ws.Range("E17").Value = "  +1.55%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0000134"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("D19").Value = "339.33"
$ws.Range("E19").Value = "  +1.59%  "

$ws.Range("D20").Value = "4.36"
$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").Value = "10.31"
$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("D22").Value = "6.40"
$ws.Range("E22").Value = "  +2.22%  "

$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").Value = "65.50"

$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("D26").Value = "0.403"
$ws.Range("E26").Value = "  -2.78%  "

$ws.Range("D27").Value = "2.723.11"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").Value = "7.05"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").Value = "0.0₃0749"
$ws.Range("E30").Value = "  -4.63%  "

$ws.Range("D32").Value = "6.25"
$ws.Range("E32").Value = "  -5.58%  "

$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("D34").Value = "18.85"
$ws.Range("E34").Value = "  +0.81%  "

$ws.Range("D35").Value = "149.84"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").Value = "4.04"
$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("E37").Value = "  -2.34%  "

$ws.Range("D38").Value = "0.869"
$ws.Range("E38").Value = "  -3.26%  "

$ws.Range("D39").Value = "0.857"
$ws.Range("E39").Value = "  +1.70%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.46"
$ws.Range("E40").Value = "  +2.32%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "36.11"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").Value = "274.09"
$ws.Range("E44").Value = "  +2.04%  "

$ws.Range("D45").Value = "0.599"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("D46").Value = "0.0959"
$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("D47").Value = "10.68"
$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("D48").Value = "18.88"
$ws.Range("E48").Value = "  -1.25%  "

$ws.Range("D49").Value = "0.0523"
$ws.Range("E49").Value = "  -1.49%  "

$ws.Range("D50").Value = "19.12"
$ws.Range("E50").Value = "  +5.29%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.978.00"
$ws.Range("E51").Value = "  -2.70%  "
